# Updates cryptos list cell values to reflect the latest scrape (GitHub Actions job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.697.34"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.546.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.39%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -5.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.75%  "
$ws.Range("E13").Value = "  +4.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.935.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.527.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.876"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.697.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0982"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.36%  "
$ws.Range("E24").Value = "  -4.92%  "
$ws.Range("E25").Value = "  -5.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.03%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.26%  "
$ws.Range("E30").Value = "  -2.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  -2.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.00%  "
$ws.Range("E36").Value = "  -5.74%  "
$ws.Range("E37").Value = "  -5.54%  "
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.64%  "
$ws.Range("E39").Value = "  -4.56%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.12"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0310"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.38%  "
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.082.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.789.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.37%  "
